$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row (row 1): swap columns so that
#   B1: bedrooms_1 -> kitchens_2
#   C1: bedrooms_2 -> bedrooms_1
#   E1: kitchens_2 -> bedrooms_2
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"

# Row 4: shift the "1" from D4 to F4
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

# Row 5: shift the "1" from B5 to E5
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1

# Row 6: shift the "1" from E6 to B6
$ws.Range("B6").Value = 1
$ws.Range("E6").Value = 0

# Row 7: shift the "1" from F7 to D7
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = 0
